$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 152, shifting rows 152:202 down to 153:203
$ws.Range("A152:R152").Insert()

# Copy formatting from the row above (151) into new row 152 so styles (e.g. date format on D) match
$ws.Range("A151:R151").Copy()
$ws.Range("A152:R152").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new row 152 values
$ws.Range("A152").Value = 8
$ws.Range("B152").Value = "Terminal La Palmera de La Serena"
$ws.Range("C152").Value = "Coquimbo"
$ws.Range("D152").Value = 44559
$ws.Range("E152").Value = 4
$ws.Range("F152").Value = 100112012
$ws.Range("G152").Value = "Espinaca"
$ws.Range("H152").Value = "Sin especificar"
$ws.Range("I152").Value = "Primera"
$ws.Range("J152").Value = 2960
$ws.Range("K152").Value = 450
$ws.Range("L152").Value = 500
$ws.Range("M152").Value = 475
$ws.Range("N152").Value = "$/atado 300 a 500 gramos"
$ws.Range("O152").Value = "Provincia del Elqu" + [char]0xED
$ws.Range("P152").Value = 950
$ws.Range("Q152").Value = 0.5
$ws.Range("R152").Value = "Hortaliza"
